$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023385068751647
$ws.Cells.Item(2, 4).Value = 1.043585436892055
$ws.Cells.Item(2, 5).Value = 1.034859835735106
$ws.Cells.Item(2, 6).Value = 1.047955900078149
$ws.Cells.Item(2, 9).Value = 1.036437067706018
$ws.Cells.Item(2, 10).Value = 1.028565585854488
$ws.Cells.Item(2, 11).Value = 1.046358887549175
$ws.Cells.Item(2, 12).Value = 1.037658099075009
$ws.Cells.Item(2, 13).Value = 1.050717091070147
$ws.Cells.Item(2, 14).Value = 1.013565718673979

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.024299847077426
$ws.Cells.Item(3, 4).Value = 1.044185885886594
$ws.Cells.Item(3, 5).Value = 1.035647395873392
$ws.Cells.Item(3, 6).Value = 1.048770402599349
$ws.Cells.Item(3, 9).Value = 1.036547826643293
$ws.Cells.Item(3, 10).Value = 1.029118999771919
$ws.Cells.Item(3, 11).Value = 1.046770952161814
$ws.Cells.Item(3, 12).Value = 1.0382549921942
$ws.Cells.Item(3, 13).Value = 1.051343535366301
$ws.Cells.Item(3, 14).Value = 1.013750277036997

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024892292105958
$ws.Cells.Item(4, 4).Value = 1.044572859244941
$ws.Cells.Item(4, 5).Value = 1.036157429088871
$ws.Cells.Item(4, 6).Value = 1.049296926171565
$ws.Cells.Item(4, 9).Value = 1.036616997160368
$ws.Cells.Item(4, 10).Value = 1.029477011433194
$ws.Cells.Item(4, 11).Value = 1.047035335644918
$ws.Cells.Item(4, 12).Value = 1.038640995729939
$ws.Cells.Item(4, 13).Value = 1.051747722940824
$ws.Cells.Item(4, 14).Value = 1.013869622249261

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025141479862231
$ws.Cells.Item(5, 4).Value = 1.044735167688048
$ws.Cells.Item(5, 5).Value = 1.036371948104555
$ws.Cells.Item(5, 6).Value = 1.049518151510465
$ws.Cells.Item(5, 9).Value = 1.036645477193856
$ws.Cells.Item(5, 10).Value = 1.029627498302021
$ws.Cells.Item(5, 11).Value = 1.047145941835028
$ws.Cells.Item(5, 12).Value = 1.038803216075573
$ws.Cells.Item(5, 13).Value = 1.051917362872507
$ws.Cells.Item(5, 14).Value = 1.013919776289624

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025183326812532
$ws.Cells.Item(6, 4).Value = 1.044762397921908
$ws.Cells.Item(6, 5).Value = 1.03640797266772
$ws.Cells.Item(6, 6).Value = 1.049555288816392
$ws.Cells.Item(6, 9).Value = 1.036650223939825
$ws.Cells.Item(6, 10).Value = 1.029652764432567
$ws.Cells.Item(6, 11).Value = 1.047164481339715
$ws.Cells.Item(6, 12).Value = 1.038830450288087
$ws.Cells.Item(6, 13).Value = 1.05194582965731
$ws.Cells.Item(6, 14).Value = 1.013928196266239

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02489562127972
$ws.Cells.Item(7, 4).Value = 1.044575029495427
$ws.Cells.Item(7, 5).Value = 1.036160295107541
$ws.Cells.Item(7, 6).Value = 1.04929988268992
$ws.Cells.Item(7, 9).Value = 1.03661738006823
$ws.Cells.Item(7, 10).Value = 1.029479022330918
$ws.Cells.Item(7, 11).Value = 1.047036815697124
$ws.Cells.Item(7, 12).Value = 1.03864316354614
$ws.Cells.Item(7, 13).Value = 1.051749990782628
$ws.Cells.Item(7, 14).Value = 1.013870292483936

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023694113525804
$ws.Cells.Item(8, 4).Value = 1.043788682491779
$ws.Cells.Item(8, 5).Value = 1.035125905503853
$ws.Cells.Item(8, 6).Value = 1.048231270186549
$ws.Cells.Item(8, 9).Value = 1.036475015676046
$ws.Cells.Item(8, 10).Value = 1.028752631450287
$ws.Cells.Item(8, 11).Value = 1.046498611370726
$ws.Cells.Item(8, 12).Value = 1.037859867913003
$ws.Cells.Item(8, 13).Value = 1.050929040675104
$ws.Cells.Item(8, 14).Value = 1.013628106614608

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021580954743713
$ws.Cells.Item(9, 4).Value = 1.042391213495345
$ws.Cells.Item(9, 5).Value = 1.033306539523106
$ws.Cells.Item(9, 6).Value = 1.046344393919148
$ws.Cells.Item(9, 9).Value = 1.036205074382693
$ws.Cells.Item(9, 10).Value = 1.027472041128935
$ws.Cells.Item(9, 11).Value = 1.045533082305144
$ws.Cells.Item(9, 12).Value = 1.036477928822416
$ws.Cells.Item(9, 13).Value = 1.049473592140324
$ws.Cells.Item(9, 14).Value = 1.013200776620715

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020174973877969
$ws.Cells.Item(10, 4).Value = 1.04145175234375
$ws.Cells.Item(10, 5).Value = 1.032095986937694
$ws.Cells.Item(10, 6).Value = 1.045084014185421
$ws.Cells.Item(10, 9).Value = 1.0360123587411
$ws.Cells.Item(10, 10).Value = 1.02661797634811
$ws.Cells.Item(10, 11).Value = 1.044877984868486
$ws.Cells.Item(10, 12).Value = 1.035555592326848
$ws.Cells.Item(10, 13).Value = 1.048497467550641
$ws.Cells.Item(10, 14).Value = 1.0129155316632

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.01956684460077
$ws.Cells.Item(11, 4).Value = 1.041043129613477
$ws.Cells.Item(11, 5).Value = 1.031572384056553
$ws.Cells.Item(11, 6).Value = 1.044537697152602
$ws.Cells.Item(11, 9).Value = 1.035925900161599
$ws.Cells.Item(11, 10).Value = 1.026248089784796
$ws.Cells.Item(11, 11).Value = 1.04459163758389
$ws.Cells.Item(11, 12).Value = 1.035155979504879
$ws.Cells.Item(11, 13).Value = 1.048073435955354
$ws.Cells.Item(11, 14).Value = 1.012791937282401

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.019341059990264
$ws.Cells.Item(12, 4).Value = 1.040891076218107
$ws.Cells.Item(12, 5).Value = 1.031377982547419
$ws.Cells.Item(12, 6).Value = 1.044334687890625
$ws.Cells.Item(12, 9).Value = 1.035893334287914
$ws.Cells.Item(12, 10).Value = 1.026110687817823
$ws.Cells.Item(12, 11).Value = 1.044484873794102
$ws.Cells.Item(12, 12).Value = 1.035007511504573
$ws.Cells.Item(12, 13).Value = 1.04791572872163
$ws.Cells.Item(12, 14).Value = 1.012746016969684

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.019389486937384
$ws.Cells.Item(13, 4).Value = 1.040923704493739
$ws.Cells.Item(13, 5).Value = 1.031419678326815
$ws.Cells.Item(13, 6).Value = 1.044378237781221
$ws.Cells.Item(13, 9).Value = 1.035900340183085
$ws.Cells.Item(13, 10).Value = 1.026140161425531
$ws.Cells.Item(13, 11).Value = 1.044507793135924
$ws.Cells.Item(13, 12).Value = 1.035039359912581
$ws.Cells.Item(13, 13).Value = 1.047949566617163
$ws.Cells.Item(13, 14).Value = 1.012755867563139

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.01954817909405
$ws.Cells.Item(14, 4).Value = 1.041030566383183
$ws.Cells.Item(14, 5).Value = 1.031556312964403
$ws.Cells.Item(14, 6).Value = 1.044520918021173
$ws.Cells.Item(14, 9).Value = 1.03592321745641
$ws.Cells.Item(14, 10).Value = 1.026236732285874
$ws.Cells.Item(14, 11).Value = 1.044582820638843
$ws.Cells.Item(14, 12).Value = 1.035143707779845
$ws.Cells.Item(14, 13).Value = 1.048060403957483
$ws.Cells.Item(14, 14).Value = 1.012788141734723

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019645968074622
$ws.Cells.Item(15, 4).Value = 1.041096371447364
$ws.Cells.Item(15, 5).Value = 1.031640509771934
$ws.Cells.Item(15, 6).Value = 1.04460881711451
$ws.Cells.Item(15, 9).Value = 1.035937253130293
$ws.Cells.Item(15, 10).Value = 1.026296231534769
$ws.Cells.Item(15, 11).Value = 1.04462899439077
$ws.Cells.Item(15, 12).Value = 1.035207995477611
$ws.Cells.Item(15, 13).Value = 1.048128667652054
$ws.Cells.Item(15, 14).Value = 1.012808025357155

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020215347616173
$ws.Cells.Item(16, 4).Value = 1.041478832923724
$ws.Cells.Item(16, 5).Value = 1.03213074897651
$ws.Cells.Item(16, 6).Value = 1.045120259728737
$ws.Cells.Item(16, 9).Value = 1.036018033361346
$ws.Cells.Item(16, 10).Value = 1.02664252309505
$ws.Cells.Item(16, 11).Value = 1.04489693236693
$ws.Cells.Item(16, 12).Value = 1.035582108473539
$ws.Cells.Item(16, 13).Value = 1.048525580562898
$ws.Cells.Item(16, 14).Value = 1.012923732531805

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.020572684412656
$ws.Cells.Item(17, 4).Value = 1.041718252558563
$ws.Cells.Item(17, 5).Value = 1.032438418002059
$ws.Cells.Item(17, 6).Value = 1.045440924500852
$ws.Cells.Item(17, 9).Value = 1.03606789922019
$ws.Cells.Item(17, 10).Value = 1.026859724388249
$ws.Cells.Item(17, 11).Value = 1.04506428503626
$ws.Cells.Item(17, 12).Value = 1.035816717671763
$ws.Cells.Item(17, 13).Value = 1.048774189844278
$ws.Cells.Item(17, 14).Value = 1.012996291051103

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020781177467085
$ws.Cells.Item(18, 4).Value = 1.041857725281028
$ws.Cells.Item(18, 5).Value = 1.032617931450227
$ws.Cells.Item(18, 6).Value = 1.045627908259124
$ws.Cells.Item(18, 9).Value = 1.036096694497611
$ws.Cells.Item(18, 10).Value = 1.026986407321082
$ws.Cells.Item(18, 11).Value = 1.045161639756379
$ws.Cells.Item(18, 12).Value = 1.035953538483523
$ws.Cells.Item(18, 13).Value = 1.048919067783507
$ws.Cells.Item(18, 14).Value = 1.013038605325615

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.020852279126579
$ws.Cells.Item(19, 4).Value = 1.041905251836535
$ws.Cells.Item(19, 5).Value = 1.03267915022924
$ws.Cells.Item(19, 6).Value = 1.045691655570476
$ws.Cells.Item(19, 9).Value = 1.036106463611919
$ws.Cells.Item(19, 10).Value = 1.027029601724771
$ws.Cells.Item(19, 11).Value = 1.045194791168381
$ws.Cells.Item(19, 12).Value = 1.036000186943489
$ws.Cells.Item(19, 13).Value = 1.048968445032292
$ws.Cells.Item(19, 14).Value = 1.013053032051251

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.020534338872093
$ws.Cells.Item(20, 4).Value = 1.041692583358513
$ws.Cells.Item(20, 5).Value = 1.032405402294539
$ws.Cells.Item(20, 6).Value = 1.045406525836359
$ws.Cells.Item(20, 9).Value = 1.036062579138926
$ws.Cells.Item(20, 10).Value = 1.026836421452042
$ws.Cells.Item(20, 11).Value = 1.045046356481717
$ws.Cells.Item(20, 12).Value = 1.035791548670243
$ws.Cells.Item(20, 13).Value = 1.048747530025211
$ws.Cells.Item(20, 14).Value = 1.012988507021221

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019501445393817
$ws.Cells.Item(21, 4).Value = 1.040999105731103
$ws.Cells.Item(21, 5).Value = 1.031516075029331
$ws.Cells.Item(21, 6).Value = 1.044478904516299
$ws.Cells.Item(21, 9).Value = 1.035916493118958
$ws.Cells.Item(21, 10).Value = 1.026208294832087
$ws.Cells.Item(21, 11).Value = 1.044560737982276
$ws.Cells.Item(21, 12).Value = 1.035112980862017
$ws.Cells.Item(21, 13).Value = 1.048027770713373
$ws.Cells.Item(21, 14).Value = 1.012778638119883

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.018852613082645
$ws.Cells.Item(22, 4).Value = 1.040561512018963
$ws.Cells.Item(22, 5).Value = 1.030957429120528
$ws.Cells.Item(22, 6).Value = 1.043895194140521
$ws.Cells.Item(22, 9).Value = 1.035822032458076
$ws.Cells.Item(22, 10).Value = 1.025813311491309
$ws.Cells.Item(22, 11).Value = 1.044253087380566
$ws.Cells.Item(22, 12).Value = 1.034686142242961
$ws.Cells.Item(22, 13).Value = 1.04757405573118
$ws.Cells.Item(22, 14).Value = 1.012646616798488

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.019196515199945
$ws.Cells.Item(23, 4).Value = 1.040793637462392
$ws.Cells.Item(23, 5).Value = 1.031253529059222
$ws.Cells.Item(23, 6).Value = 1.044204674708619
$ws.Cells.Item(23, 9).Value = 1.035872354923127
$ws.Cells.Item(23, 10).Value = 1.02602270452714
$ws.Cells.Item(23, 11).Value = 1.044416398366549
$ws.Cells.Item(23, 12).Value = 1.034912435626128
$ws.Cells.Item(23, 13).Value = 1.047814689292678
$ws.Cells.Item(23, 14).Value = 1.012716610158148

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02055166537193
$ws.Cells.Item(24, 4).Value = 1.041704182711252
$ws.Cells.Item(24, 5).Value = 1.032420320500869
$ws.Cells.Item(24, 6).Value = 1.045422069280963
$ws.Cells.Item(24, 9).Value = 1.036064983952738
$ws.Cells.Item(24, 10).Value = 1.02684695106727
$ws.Cells.Item(24, 11).Value = 1.045054458424997
$ws.Cells.Item(24, 12).Value = 1.035802921529086
$ws.Cells.Item(24, 13).Value = 1.048759576856793
$ws.Cells.Item(24, 14).Value = 1.012992024313652

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022126769841063
$ws.Cells.Item(25, 4).Value = 1.042753880004794
$ws.Cells.Item(25, 5).Value = 1.033776480787627
$ws.Cells.Item(25, 6).Value = 1.046832639545364
$ws.Cells.Item(25, 9).Value = 1.036277113963760
$ws.Cells.Item(25, 10).Value = 1.027803168593636
$ws.Cells.Item(25, 11).Value = 1.045784714398231
$ws.Cells.Item(25, 12).Value = 1.036835383025759
$ws.Cells.Item(25, 13).Value = 1.049850895314679
$ws.Cells.Item(25, 14).Value = 1.013311316516233
